$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MainSheet")

# --- Update row 3 and row 4 values (A3:E4) ---
# Row 3 becomes what Row 4 previously held (Run, webApp, Web, Windows_Chrome, None)
# Row 4 keeps its B:E values, but its A value swaps to what Row 3 previously held (No)
$ws.Range("A3").Value = "Run"
$ws.Range("B3").Value = "webApp"
$ws.Range("C3").Value = "Web"
$ws.Range("D3").Value = "Windows_Chrome"
$ws.Range("E3").Value = "None"

$ws.Range("A4").Value = "No"
$ws.Range("B4").Value = "webApp"
$ws.Range("C4").Value = "Web"
$ws.Range("D4").Value = "Windows_Chrome"
$ws.Range("E4").Value = "None"

# --- Update the selection to B3 (active cell) ---
$ws.Range("B3").Select()

# --- Data validation changes ---
# Narrow the existing "webApp,nativeApp,sanityTesting" list validation to B3:B4 only
$ws.Range("B2:B4").Validation.Delete()
$ws.Range("B3:B4").Validation.Add(3, 1, 1, "webApp,nativeApp,sanityTesting")
$ws.Range("B3:B4").Validation.InputMessage = ""
$ws.Range("B3:B4").Validation.ErrorMessage = ""
$ws.Range("B3:B4").Validation.ShowInput = $true
$ws.Range("B3:B4").Validation.ShowError = $true

# Add a new validation list on B2 only with a distinct option set
$ws.Range("B2").Validation.Add(3, 1, 1, "webApp,NativeApp,sanityTesting")
$ws.Range("B2").Validation.InputMessage = ""
$ws.Range("B2").Validation.ErrorMessage = ""
$ws.Range("B2").Validation.ShowInput = $true
$ws.Range("B2").Validation.ShowError = $true
